$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to remain
# text (matching the source data, which stores prices as literal strings),
# otherwise Excel auto-converts them to numeric values.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D16",
    "D22",
    "D23",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D49",
    "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the new cell values from the commit diff.
$newValues = [ordered]@{
    "D2" = "42.535.26"
    "E2" = "  -2.87%  "
    "D3" = "2.230.10"
    "E3" = "  -1.95%  "
    "E4" = "  +0.05%  "
    "D5" = "112.65"
    "E5" = "  -7.42%  "
    "D6" = "297.20"
    "E6" = "  +11.20%  "
    "D7" = "0.622"
    "E7" = "  -4.49%  "
    "E8" = "  -0.41%  "
    "D9" = "0.608"
    "E9" = "  -3.19%  "
    "D10" = "44.49"
    "E10" = "  -8.38%  "
    "D11" = "0.0917"
    "E11" = "  -3.31%  "
    "D12" = "54.46"
    "E12" = "  +0.33%  "
    "D13" = "8.87"
    "E13" = "  -4.32%  "
    "E14" = "  +10.61%  "
    "E15" = "  -2.96%  "
    "D16" = "15.13"
    "E16" = "  -3.39%  "
    "D17" = "2.562.87"
    "E17" = "  -2.08%  "
    "D18" = "2.242.49"
    "E18" = "  -1.52%  "
    "D19" = "42.526.55"
    "E19" = "  -2.70%  "
    "E20" = "  +5.25%  "
    "E21" = "  -4.08%  "
    "D22" = "73.08"
    "D23" = "3.52"
    "E23" = "  +21.29%  "
    "E24" = "  -2.06%  "
    "D25" = "230.02"
    "E25" = "  -2.49%  "
    "D26" = "9.25"
    "E26" = "  -3.62%  "
    "E27" = "  -1.58%  "
    "D28" = "11.65"
    "E28" = "  -2.80%  "
    "B29" = "InjectiveProtocol"
    "C29" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D29" = "38.78"
    "E29" = "  -10.51%  "
    "B30" = "Toncoin"
    "C30" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "D30" = "2.23"
    "E30" = "  -1.19%  "
    "E31" = "  -5.70%  "
    "D32" = "173.28"
    "E32" = "  -0.44%  "
    "D33" = "21.07"
    "E33" = "  -2.92%  "
    "D34" = "0.0898"
    "E34" = "  -3.25%  "
    "D35" = "5.78"
    "E35" = "  -0.28%  "
    "D36" = "5.17"
    "E36" = "  +11.69%  "
    "D37" = "4.32"
    "E37" = "  +0.96%  "
    "E38" = "  -3.95%  "
    "D39" = "0.0377"
    "E39" = "  -1.95%  "
    "E40" = "  -4.28%  "
    "D41" = "2.42"
    "E41" = "  -5.34%  "
    "D42" = "72.15"
    "E42" = "  -2.75%  "
    "D43" = "0.237"
    "E43" = "  -1.17%  "
    "B44" = "Celestia"
    "C44" = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
    "D44" = "12.74"
    "E44" = "  -7.31%  "
    "B45" = "FirstDigitalUSD"
    "C45" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D45" = "1.00"
    "E45" = "  -0.12%  "
    "E46" = "  -4.70%  "
    "D47" = "5.44"
    "E47" = "  -7.96%  "
    "E48" = "  +3.41%  "
    "D49" = "103.26"
    "E49" = "  -0.17%  "
    "E50" = "  +0.85%  "
    "D51" = "1.66"
    "E51" = "  +7.55%  "
}
foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
